# "Stun now cancels abilities"
# Rebalance castTime for Heal, Poke, Nuke, Stun; buff Nuke power; underline
# Stun's castTime cell to flag the new cancel-on-stun behavior.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Heal castTime: 30 -> 3
$ws.Range("D2").Value = 3

# Poke castTime: 50 -> 5
$ws.Range("D3").Value = 5

# Nuke power: 2 -> 5
$ws.Range("B5").Value = 5

# Nuke castTime: 20 -> 2
$ws.Range("D5").Value = 2

# Stun castTime: 50 -> 5, underlined to highlight the new cancel behavior
$ws.Range("D6").Value = 5
$ws.Range("D6").Font.Underline = [int]2  # xlUnderlineStyleSingle

# Reflect the editor's final selection
$ws.Range("B5").Select()
